$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Shared-string edits -------------------------------------------------
# "Residential" (header J1) becomes "Residentia"
$ws.Range("J1").Value = "Residentia"

# New column K: header "houseex" plus per-province change-in-expenditure data
$ws.Range("K1").Value = "houseex"

$kValues = @{
    2  = 0
    3  = 0
    4  = 0
    5  = 0
    6  = 0
    7  = -4.9292343457456536
    8  = 0
    9  = 0
    10 = -4.8203515127295331
    11 = 0
    12 = 0
    13 = 0
    14 = -7.5196914542577344
    15 = 0
    16 = 0
    17 = 0
    18 = 0
    19 = 0
    20 = 0
    21 = 0
    22 = 0
    23 = 0
    24 = 0
    25 = 0
    26 = 0
    27 = -6.181305037989234
    28 = 0
    29 = 0
    30 = -6.3433217044006769
    31 = 0
    32 = 0
    33 = 0
    34 = 0
    35 = 0
}

foreach ($row in $kValues.Keys) {
    $ws.Range("K$row").Value = $kValues[$row]
}

# New column K width (auto-fit-like, as close as the host allows)
$ws.Columns.Item(11).ColumnWidth = 14

# Restore the view/selection state as closely as the host allows
$win = $excel.ActiveWindow
$win.ScrollRow = 7
$win.ScrollColumn = 1
$ws.Range("K39").Select()
